# Update "想去人数" (column F) counts on the "展览" sheet and the mirrored
# "全部类型" aggregate sheet, matching the regenerated gh-pages data output.

$wb = $excel.ActiveWorkbook

$exhibitionUpdates = @{
    4  = 124
    5  = 513
    6  = 725
    7  = 1419
    9  = 80
    10 = 126
    11 = 6062
    12 = 64
    13 = 396
    14 = 106
    15 = 4872
    18 = 1159
    19 = 48
    20 = 350
    21 = 51
    23 = 276
    24 = 22
    25 = 3347
    26 = 138
}

$allTypesUpdates = @{
    4  = 124
    6  = 513
    7  = 725
    8  = 1419
    10 = 80
    11 = 126
    12 = 6062
    13 = 64
    14 = 396
    15 = 106
    16 = 4872
    19 = 1159
    20 = 48
    21 = 350
    22 = 51
    24 = 276
    25 = 22
    26 = 3347
    28 = 138
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
